# Apply the "Updated symbol list" commit: refresh Price (column D) and a
# couple of Volume(1h) label (column E) cells on the active worksheet.
#
# The Price column stores numeric-looking values as TEXT (inline strings in
# the OOXML). A plain numeric-looking string assigned to .Value gets
# auto-converted by Excel into a real number, which would change the cell's
# stored type. To keep these as text (matching the original authoring),
# each value is prefixed with a leading apostrophe, which is Excel's
# standard "force text" entry convention and is stripped from the stored
# text value itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'236.82"
$ws.Range("D3").Value = "'21.65"
$ws.Range("D4").Value = "'5.363"
$ws.Range("D5").Value = "'0.05569"
$ws.Range("D6").Value = "'3.367"
$ws.Range("D8").Value = "'0.7987"
$ws.Range("D9").Value = "'1.031"
$ws.Range("D11").Value = "'0.07306"
$ws.Range("D12").Value = "'0.03207"
$ws.Range("D13").Value = "'0.02914"
$ws.Range("D14").Value = "'0.09269"
$ws.Range("D15").Value = "'0.001658"
$ws.Range("D16").Value = "'3.254"
$ws.Range("D18").Value = "'0.0005708"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006258"
$ws.Range("D20").Value = "'0.005071"
$ws.Range("D23").Value = "'0.0004180"
$ws.Range("D24").Value = "'3.952"
$ws.Range("D27").Value = "'0.1295"
$ws.Range("D40").Value = "'0.04128"
$ws.Range("D41").Value = "'0.007016"
$ws.Range("D42").Value = "'0.003497"
$ws.Range("D43").Value = "'0.1036"
$ws.Range("D44").Value = "'0.009517"
$ws.Range("D45").Value = "'0.00005439"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D47").Value = "'0.6797"
$ws.Range("D48").Value = "'0.03229"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").Value = "'0.00002098"
